# Auto-generated Excel COM-interop script
# Applies cached market-price value updates to the per-job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as captured by the
# scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 49998
$ws.Range("J3").Value = 49998
$ws.Range("L3").Value = 49998
$ws.Range("N3").Value = -50226
$ws.Range("H9").Value = 72.6
$ws.Range("I9").Value = 77
$ws.Range("K9").Value = 77
$ws.Range("M9").Value = 92
$ws.Range("H11").Value = 56.8
$ws.Range("I11").Value = 56.8
$ws.Range("K11").Value = 56.8
$ws.Range("M11").Value = 83.2
$ws.Range("H17").Value = 1689.1052
$ws.Range("I17").Value = 999.9167
$ws.Range("J17").Value = 2870.5715
$ws.Range("K17").Value = 2999.7501
$ws.Range("L17").Value = 8611.7145
$ws.Range("M17").Value = -2831.7501
$ws.Range("N17").Value = -8947.7145
$ws.Range("H19").Value = 2067.6
$ws.Range("I19").Value = 2328.125
$ws.Range("J19").Value = 1025.5
$ws.Range("K19").Value = 2328.125
$ws.Range("L19").Value = 1025.5
$ws.Range("M19").Value = -2153.125
$ws.Range("N19").Value = -1375.5
$ws.Range("H62").Value = 4900
$ws.Range("I62").Value = 3800
$ws.Range("K62").Value = 3800
$ws.Range("M62").Value = -3176
$ws.Range("H65").Value = 4900
$ws.Range("I65").Value = 3800
$ws.Range("K65").Value = 19000
$ws.Range("M65").Value = -15880
$ws.Range("H70").Value = 4874.5
$ws.Range("J70").Value = 5749
$ws.Range("L70").Value = 17247
$ws.Range("N70").Value = -17787
$ws.Range("H73").Value = 4874.5
$ws.Range("J73").Value = 5749
$ws.Range("L73").Value = 17247
$ws.Range("N73").Value = -19119
$ws.Range("H93").Value = 14471.5
$ws.Range("I93").Value = 13999
$ws.Range("J93").Value = 14629
$ws.Range("K93").Value = 13999
$ws.Range("L93").Value = 14629
$ws.Range("M93").Value = -11503
$ws.Range("N93").Value = -19621
$ws.Range("H102").Value = 49998
$ws.Range("J102").Value = 49998
$ws.Range("L102").Value = 49998
$ws.Range("N102").Value = -56488
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 6999.6665
$ws.Range("I132").Value = 6999.6665
$ws.Range("K132").Value = 20998.9995
$ws.Range("M132").Value = -18468.9995
$ws.Range("H138").Value = 3132.318
$ws.Range("J138").Value = 3158.1
$ws.Range("L138").Value = 9474.3
$ws.Range("N138").Value = -19754.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100.666664
$ws.Range("I5").Value = 94.333336
$ws.Range("J5").Value = 113.333336
$ws.Range("K5").Value = 94.333336
$ws.Range("L5").Value = 113.333336
$ws.Range("M5").Value = 17.666664
$ws.Range("N5").Value = -337.333336
$ws.Range("H6").Value = 1800.4
$ws.Range("I6").Value = 1800.4
$ws.Range("K6").Value = 1800.4
$ws.Range("M6").Value = -1627.4
$ws.Range("H42").Value = 2028
$ws.Range("I42").Value = 2028
$ws.Range("K42").Value = 2028
$ws.Range("M42").Value = -1542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100.666664
$ws.Range("I4").Value = 94.333336
$ws.Range("J4").Value = 113.333336
$ws.Range("K4").Value = 94.333336
$ws.Range("L4").Value = 113.333336
$ws.Range("M4").Value = 20.666664
$ws.Range("N4").Value = -343.333336
$ws.Range("H22").Value = 268.25
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.066666
$ws.Range("I7").Value = 52.833332
$ws.Range("K7").Value = 52.833332
$ws.Range("M7").Value = 60.166668
$ws.Range("H25").Value = 500
$ws.Range("J25").Value = 500
$ws.Range("L25").Value = 500
$ws.Range("N25").Value = -848
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H60").Value = 4400
$ws.Range("I60").Value = 4400
$ws.Range("K60").Value = 4400
$ws.Range("M60").Value = -3889
$ws.Range("H86").Value = 12500
$ws.Range("J86").Value = 12500
$ws.Range("L86").Value = 12500
$ws.Range("N86").Value = -14746
$ws.Range("H89").Value = 12500
$ws.Range("J89").Value = 12500
$ws.Range("L89").Value = 62500
$ws.Range("N89").Value = -73732
$ws.Range("H107").Value = 2349.4
$ws.Range("I107").Value = 1375
$ws.Range("K107").Value = 1375
$ws.Range("M107").Value = 545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2499
$ws.Range("I63").Value = 2499
$ws.Range("K63").Value = 7497
$ws.Range("M63").Value = -6748
$ws.Range("H66").Value = 2499
$ws.Range("I66").Value = 2499
$ws.Range("K66").Value = 22491
$ws.Range("M66").Value = -18747
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 15000
$ws.Range("M70").Value = -14685
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 15000
$ws.Range("M73").Value = -13908
$ws.Range("H103").Value = 4113.7144
$ws.Range("J103").Value = 4959.4
$ws.Range("L103").Value = 14878.2
$ws.Range("N103").Value = -16636.2
$ws.Range("H109").Value = 2969.5
$ws.Range("I109").Value = 2969.5
$ws.Range("K109").Value = 8908.5
$ws.Range("M109").Value = -7868.5
$ws.Range("H129").Value = 562.7143
$ws.Range("I129").Value = 599.8
$ws.Range("J129").Value = 470
$ws.Range("K129").Value = 1799.4
$ws.Range("L129").Value = 1410
$ws.Range("M129").Value = 3200.6
$ws.Range("N129").Value = -11410
$ws.Range("H131").Value = 3998.2144
$ws.Range("I131").Value = 6407.1665
$ws.Range("J131").Value = 2191.5
$ws.Range("K131").Value = 19221.4995
$ws.Range("L131").Value = 6574.5
$ws.Range("M131").Value = -14181.4995
$ws.Range("N131").Value = -16654.5
$ws.Range("H139").Value = 2639.7144
$ws.Range("I139").Value = 2639.7144
$ws.Range("K139").Value = 7919.1432
$ws.Range("M139").Value = -2779.1432
$ws.Range("H140").Value = 1936.75
$ws.Range("I140").Value = 1936.75
$ws.Range("K140").Value = 5810.25
$ws.Range("M140").Value = -630.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.888885
$ws.Range("I2").Value = 46.6
$ws.Range("K2").Value = 46.6
$ws.Range("M2").Value = 66.4
$ws.Range("H4").Value = 3000
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1625
$ws.Range("I22").Value = 1625
$ws.Range("J22").Value = 1625
$ws.Range("K22").Value = 1625
$ws.Range("L22").Value = 1625
$ws.Range("M22").Value = -1330
$ws.Range("N22").Value = -2215
$ws.Range("H27").Value = 1625
$ws.Range("I27").Value = 1625
$ws.Range("J27").Value = 1625
$ws.Range("K27").Value = 1625
$ws.Range("L27").Value = 1625
$ws.Range("M27").Value = -1518
$ws.Range("N27").Value = -1839
$ws.Range("H46").Value = 6508
$ws.Range("I46").Value = 7450
$ws.Range("J46").Value = 4624
$ws.Range("K46").Value = 7450
$ws.Range("L46").Value = 4624
$ws.Range("M46").Value = -7262
$ws.Range("N46").Value = -5000
$ws.Range("H132").Value = 9000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 32000
$ws.Range("I134").Value = 32000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 32000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -26930
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 17686.334
$ws.Range("I51").Value = 17686.334
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 17686.334
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -17176.334
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 4999
$ws.Range("J62").Value = 4999
$ws.Range("L62").Value = 4999
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 4999
$ws.Range("J65").Value = 4999
$ws.Range("L65").Value = 24995
$ws.Range("N65").Value = -31235
$ws.Range("H132").Value = 7050
$ws.Range("I132").Value = 7050
$ws.Range("K132").Value = 21150
$ws.Range("M132").Value = -18620
